$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.930.92"
$ws.Range("E2").Value = "  +0.22%  "

# Row 3
$ws.Range("D3").Value = "2.824.27"
$ws.Range("E3").Value = "  +1.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.87"
$ws.Range("E5").Value = "  +6.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.56"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("E7").Value = "  +4.45%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +5.38%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.57"
$ws.Range("E10").Value = "  -1.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  -0.08%  "

# Row 12
$ws.Range("E12").Value = "  +1.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.93"
$ws.Range("E13").Value = "  -1.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.76"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15
$ws.Range("D15").Value = "3.267.76"
$ws.Range("E15").Value = "  +1.76%  "

# Row 16
$ws.Range("D16").Value = "2.818.68"
$ws.Range("E16").Value = "  +2.02%  "

# Row 17
$ws.Range("E17").Value = "  -0.16%  "

# Row 18
$ws.Range("D18").Value = "51.847.31"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("E19").Value = "  +9.03%  "

# Row 20
$ws.Range("E20").Value = "  -2.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.45"
$ws.Range("E21").Value = "  -0.92%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0987"
$ws.Range("E22").Value = "  +1.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.60"
$ws.Range("E23").Value = "  -3.63%  "

# Row 24
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +3.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.80"
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("E28").Value = "  +0.94%  "

# Row 29
$ws.Range("E29").Value = "  +0.84%  "

# Row 30
$ws.Range("E30").Value = "  -2.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "50.77"
$ws.Range("E31").Value = "  +1.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.03"
$ws.Range("E32").Value = "  -3.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0454"
$ws.Range("E33").Value = "  +25.92%  "

# Row 34
$ws.Range("E34").Value = "  +4.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  +5.95%  "

# Row 36
$ws.Range("E36").Value = "  +0.34%  "

# Row 37
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("E38").Value = "  -1.98%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.22"
$ws.Range("E39").Value = "  -1.83%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.24"
$ws.Range("E40").Value = "  -5.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.98"
$ws.Range("E41").Value = "  +2.55%  "

# Row 42
$ws.Range("E42").Value = "  +2.07%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.98"
$ws.Range("E44").Value = "  -1.27%  "

# Row 45
$ws.Range("E45").Value = "  -1.57%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.38"
$ws.Range("E46").Value = "  +1.46%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.094.39"
$ws.Range("E47").Value = "  +0.15%  "

# Row 48
$ws.Range("E48").Value = "  +1.50%  "

# Row 49
$ws.Range("E49").Value = "  +3.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.930"
$ws.Range("E50").Value = "  +5.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.88"
$ws.Range("E51").Value = "  +0.03%  "
